$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new date
$ws.Name = "Through 2022-08-28"

# Update header label for the "through" date
$ws.Range("I1").Value = "2022 (through 08-28)"

# Update September (row 9) value for 2022 column (I)
$ws.Range("I9").Value = 149

# Update Total (row 14) value for 2022 column (I)
$ws.Range("I14").Value = 1120
